$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10799.2
$ws.Range("I18").Value = 10799.2
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 10799.2
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents() | Out-Null
$ws.Range("N18").Value = -10515.2

$ws.Range("H98").Value = 20325.697
$ws.Range("I98").Value = 20814.428
$ws.Range("J98").Value = 17588.8
$ws.Range("K98").Value = 20814.428
$ws.Range("L98").Value = 17588.8
$ws.Range("M98").Value = -19316.428
$ws.Range("N98").Value = -20584.8

$ws.Range("H106").Value = 6865117.5
$ws.Range("I106").Value = 11228148
$ws.Range("J106").Value = 8926.429
$ws.Range("K106").Value = 11228148
$ws.Range("L106").Value = 8926.429
$ws.Range("M106").Value = -11227517
$ws.Range("N106").Value = -10188.429

$ws.Range("H107").Value = 4795.4707
$ws.Range("I107").Value = 4760.3105
$ws.Range("J107").Value = 4999.4
$ws.Range("K107").Value = 4760.3105
$ws.Range("L107").Value = 4999.4
$ws.Range("M107").Value = -2840.3105
$ws.Range("N107").Value = -8839.4

$ws.Range("H108").Value = 98000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 98000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 98000
$ws.Range("N108").Value = -105680

$ws.Range("H122").Value = 20325.697
$ws.Range("I122").Value = 20814.428
$ws.Range("J122").Value = 17588.8
$ws.Range("K122").Value = 62443.284
$ws.Range("L122").Value = 52766.39999999999
$ws.Range("M122").Value = -59993.284
$ws.Range("N122").Value = -57666.39999999999

$ws.Range("H129").Value = 1056.4546
$ws.Range("I129").Value = 912.1
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 2736.3
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = 2263.7
$ws.Range("N129").Value = -17500

$ws.Range("H132").Value = 6856.095
$ws.Range("I132").Value = 7677.0625
$ws.Range("J132").Value = 4229
$ws.Range("K132").Value = 23031.1875
$ws.Range("L132").Value = 12687
$ws.Range("M132").Value = -20501.1875
$ws.Range("N132").Value = -17747

$ws.Range("H138").Value = 2057.587
$ws.Range("I138").Value = 1092
$ws.Range("J138").Value = 3868.0625
$ws.Range("K138").Value = 3276
$ws.Range("L138").Value = 11604.1875
$ws.Range("M138").Value = 1864
$ws.Range("N138").Value = -21884.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7583.148
$ws.Range("I32").Value = 7839.8
$ws.Range("J32").Value = 4375
$ws.Range("K32").Value = 7839.8
$ws.Range("L32").Value = 4375
$ws.Range("M32").Value = -7552.8
$ws.Range("N32").Value = -4949

$ws.Range("H44").Value = 79024.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 79024.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 79024.5
$ws.Range("N44").Value = -80000.5

$ws.Range("H61").Value = 4585.75
$ws.Range("I61").Value = 4829.5527
$ws.Range("J61").Value = 3041.6667
$ws.Range("K61").Value = 4829.5527
$ws.Range("L61").Value = 3041.6667
$ws.Range("M61").Value = -4617.5527
$ws.Range("N61").Value = -3465.6667

$ws.Range("H136").Value = 4585.75
$ws.Range("I136").Value = 4829.5527
$ws.Range("J136").Value = 3041.6667
$ws.Range("K136").Value = 14488.6581
$ws.Range("L136").Value = 9125.000100000001
$ws.Range("M136").Value = -11938.6581
$ws.Range("N136").Value = -14225.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents() | Out-Null
$ws.Range("N132").Value = 0

$ws.Range("H134").Value = 15538.611
$ws.Range("I134").Value = 17081
$ws.Range("J134").Value = 3199.5
$ws.Range("K134").Value = 51243
$ws.Range("L134").Value = 9598.5
$ws.Range("M134").Value = -48708
$ws.Range("N134").Value = -14668.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 792.9091
$ws.Range("I22").Value = 427.06668
$ws.Range("J22").Value = 1576.8572
$ws.Range("K22").Value = 427.06668
$ws.Range("L22").Value = 1576.8572
$ws.Range("M22").Value = -77.06668000000002
$ws.Range("N22").Value = -2276.8572

$ws.Range("H31").Value = 3067.75
$ws.Range("I31").Value = 3106.122
$ws.Range("J31").Value = 2984.9473
$ws.Range("K31").Value = 3106.122
$ws.Range("L31").Value = 2984.9473
$ws.Range("M31").Value = -2811.122
$ws.Range("N31").Value = -3574.9473

$ws.Range("H34").Value = 3067.75
$ws.Range("I34").Value = 3106.122
$ws.Range("J34").Value = 2984.9473
$ws.Range("K34").Value = 3106.122
$ws.Range("L34").Value = 2984.9473
$ws.Range("M34").Value = -2904.122
$ws.Range("N34").Value = -3388.9473

$ws.Range("H74").Value = 49997.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 49997.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 49997.5
$ws.Range("N74").Value = -51745.5

$ws.Range("H77").Value = 49997.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 49997.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 149992.5
$ws.Range("N77").Value = -158728.5

$ws.Range("H99").Value = 317502.88
$ws.Range("I99").Value = 458419.9
$ws.Range("J99").Value = 7485.4
$ws.Range("K99").Value = 458419.9
$ws.Range("L99").Value = 7485.4
$ws.Range("M99").Value = -456921.9
$ws.Range("N99").Value = -10481.4

$ws.Range("H126").Value = 317502.88
$ws.Range("I126").Value = 458419.9
$ws.Range("J126").Value = 7485.4
$ws.Range("K126").Value = 1375259.7
$ws.Range("L126").Value = 22456.2
$ws.Range("M126").Value = -1372789.7
$ws.Range("N126").Value = -27396.2

$ws.Range("H132").Value = 37459.547
$ws.Range("I132").Value = 1254.125
$ws.Range("J132").Value = 134007.33
$ws.Range("K132").Value = 3762.375
$ws.Range("L132").Value = 402021.99
$ws.Range("M132").Value = -1232.375
$ws.Range("N132").Value = -407081.99

$ws.Range("H134").Value = 1738.5536
$ws.Range("I134").Value = 1462.8959
$ws.Range("J134").Value = 3392.5
$ws.Range("K134").Value = 4388.6877
$ws.Range("L134").Value = 10177.5
$ws.Range("M134").Value = -1853.6877
$ws.Range("N134").Value = -15247.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10421.1
$ws.Range("I87").Value = 6635.143
$ws.Range("J87").Value = 19255
$ws.Range("K87").Value = 19905.429
$ws.Range("L87").Value = 57765
$ws.Range("M87").Value = -18657.429
$ws.Range("N87").Value = -60261

$ws.Range("H90").Value = 10421.1
$ws.Range("I90").Value = 6635.143
$ws.Range("J90").Value = 19255
$ws.Range("K90").Value = 59716.287
$ws.Range("L90").Value = 173295
$ws.Range("M90").Value = -53476.287
$ws.Range("N90").Value = -185775

$ws.Range("H97").Value = 35024.89
$ws.Range("I97").Value = 55700.637
$ws.Range("J97").Value = 2534.4285
$ws.Range("K97").Value = 167101.911
$ws.Range("L97").Value = 7603.2855
$ws.Range("M97").Value = -166605.911
$ws.Range("N97").Value = -8595.2855

$ws.Range("H131").Value = 4528.407
$ws.Range("I131").Value = 15553.2
$ws.Range("J131").Value = 2022.7727
$ws.Range("K131").Value = 46659.60000000001
$ws.Range("L131").Value = 6068.3181
$ws.Range("M131").Value = -41619.60000000001
$ws.Range("N131").Value = -16148.3181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8670.454
$ws.Range("I80").Value = 10879.583
$ws.Range("J80").Value = 6019.5
$ws.Range("K80").Value = 10879.583
$ws.Range("L80").Value = 6019.5
$ws.Range("M80").Value = -9881.583000000001
$ws.Range("N80").Value = -8015.5

$ws.Range("H83").Value = 8670.454
$ws.Range("I83").Value = 10879.583
$ws.Range("J83").Value = 6019.5
$ws.Range("K83").Value = 54397.915
$ws.Range("L83").Value = 30097.5
$ws.Range("M83").Value = -49405.915
$ws.Range("N83").Value = -40081.5

$ws.Range("H102").Value = 6297.41
$ws.Range("I102").Value = 6634.4707
$ws.Range("J102").Value = 4005.4
$ws.Range("K102").Value = 6634.4707
$ws.Range("L102").Value = 4005.4
$ws.Range("M102").Value = -5012.4707
$ws.Range("N102").Value = -7249.4

$ws.Range("H132").Value = 3871.608
$ws.Range("I132").Value = 3390.7073
$ws.Range("J132").Value = 5843.3
$ws.Range("K132").Value = 10172.1219
$ws.Range("L132").Value = 17529.9
$ws.Range("M132").Value = -7642.1219
$ws.Range("N132").Value = -22589.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22931.88
$ws.Range("I7").Value = 61685.715
$ws.Range("J7").Value = 7860.9443
$ws.Range("K7").Value = 61685.715
$ws.Range("L7").Value = 7860.9443
$ws.Range("M7").Value = -61573.715
$ws.Range("N7").Value = -8084.9443

$ws.Range("H22").Value = 16424.715
$ws.Range("I22").Value = 36817.332
$ws.Range("J22").Value = 1130.25
$ws.Range("K22").Value = 36817.332
$ws.Range("L22").Value = 1130.25
$ws.Range("M22").Value = -36522.332
$ws.Range("N22").Value = -1720.25

$ws.Range("H27").Value = 16424.715
$ws.Range("I27").Value = 36817.332
$ws.Range("J27").Value = 1130.25
$ws.Range("K27").Value = 36817.332
$ws.Range("L27").Value = 1130.25
$ws.Range("M27").Value = -36710.332
$ws.Range("N27").Value = -1344.25

$ws.Range("H40").Value = 14231.932
$ws.Range("I40").Value = 15414.866
$ws.Range("J40").Value = 11697.071
$ws.Range("K40").Value = 15414.866
$ws.Range("L40").Value = 11697.071
$ws.Range("M40").Value = -15278.866
$ws.Range("N40").Value = -11969.071

$ws.Range("H55").Value = 1554.909
$ws.Range("I55").Value = 222
$ws.Range("J55").Value = 3887.5
$ws.Range("K55").Value = 222
$ws.Range("L55").Value = 3887.5
$ws.Range("M55").Value = -49
$ws.Range("N55").Value = -4233.5

$ws.Range("H93").Value = 5765.3105
$ws.Range("I93").Value = 6791.522
$ws.Range("J93").Value = 1831.5
$ws.Range("K93").Value = 6791.522
$ws.Range("L93").Value = 1831.5
$ws.Range("M93").Value = -5543.522
$ws.Range("N93").Value = -4327.5

$ws.Range("H122").Value = 12925
$ws.Range("I122").Value = 14200
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 42600
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -40150
$ws.Range("N122").Value = -16900

$ws.Range("H126").Value = 22931.88
$ws.Range("I126").Value = 61685.715
$ws.Range("J126").Value = 7860.9443
$ws.Range("K126").Value = 185057.145
$ws.Range("L126").Value = 23582.8329
$ws.Range("M126").Value = -182587.145
$ws.Range("N126").Value = -28522.8329

$ws.Range("H132").Value = 278823.72
$ws.Range("I132").Value = 348609.1
$ws.Range("J132").Value = 6026.273
$ws.Range("K132").Value = 1045827.3
$ws.Range("L132").Value = 18078.819
$ws.Range("M132").Value = -1043297.3
$ws.Range("N132").Value = -23138.819

$ws.Range("H136").Value = 4696.0464
$ws.Range("I136").Value = 2070.0417
$ws.Range("J136").Value = 8013.1055
$ws.Range("K136").Value = 6210.125100000001
$ws.Range("L136").Value = 24039.3165
$ws.Range("M136").Value = -3660.125100000001
$ws.Range("N136").Value = -29139.3165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 5912.6665
$ws.Range("I29").Value = 5912.6665
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5912.6665
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents() | Out-Null
$ws.Range("N29").Value = -5622.6665

$ws.Range("H80").Value = 20000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 20000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -21996

$ws.Range("H83").Value = 20000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 20000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -69984

$ws.Range("H96").Value = 5557625
$ws.Range("I96").Value = 11112767
$ws.Range("J96").Value = 2483.6667
$ws.Range("K96").Value = 11112767
$ws.Range("L96").Value = 2483.6667
$ws.Range("M96").Value = -11111394
$ws.Range("N96").Value = -5229.6667

$ws.Range("H132").Value = 13416.2
$ws.Range("I132").Value = 14871.544
$ws.Range("J132").Value = 5977.778
$ws.Range("K132").Value = 44614.632
$ws.Range("L132").Value = 17933.334
$ws.Range("M132").Value = -42084.632
$ws.Range("N132").Value = -22993.334

$ws.Range("H136").Value = 230808.86
$ws.Range("I136").Value = 275553.38
$ws.Range("J136").Value = 3018.6365
$ws.Range("K136").Value = 826660.14
$ws.Range("L136").Value = 9055.9095
$ws.Range("M136").Value = -824110.14
$ws.Range("N136").Value = -14155.9095
